# Auto-generated script applying cached-price/profit updates from the
# scheduled market-data refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 46.666668
$ws.Range("I5").Value = 48.4
$ws.Range("J5").Value = 38
$ws.Range("K5").Value = 48.4
$ws.Range("L5").Value = 38
$ws.Range("M5").Value = 66.59999999999999
$ws.Range("N5").Value = -268
$ws.Range("H19").Value = 1505.037
$ws.Range("J19").Value = 1068
$ws.Range("L19").Value = 1068
$ws.Range("N19").Value = -1418
$ws.Range("H62").Value = 7402.4443
$ws.Range("I62").Value = 7077.75
$ws.Range("K62").Value = 7077.75
$ws.Range("M62").Value = -6453.75
$ws.Range("H64").Value = 7375.125
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 8800.200000000001
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 8800.200000000001
$ws.Range("M64").Value = -4752
$ws.Range("N64").Value = -9296.200000000001
$ws.Range("H65").Value = 7402.4443
$ws.Range("I65").Value = 7077.75
$ws.Range("K65").Value = 35388.75
$ws.Range("M65").Value = -32268.75
$ws.Range("H67").Value = 7375.125
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 8800.200000000001
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 8800.200000000001
$ws.Range("M67").Value = -4142
$ws.Range("N67").Value = -10516.2
$ws.Range("H70").Value = 22731228
$ws.Range("I70").Value = 792.3333
$ws.Range("K70").Value = 2376.9999
$ws.Range("M70").Value = -2106.9999
$ws.Range("H73").Value = 22731228
$ws.Range("I73").Value = 792.3333
$ws.Range("K73").Value = 2376.9999
$ws.Range("M73").Value = -1440.9999
$ws.Range("H98").Value = 3343.4285
$ws.Range("I98").Value = 2849.75
$ws.Range("K98").Value = 2849.75
$ws.Range("M98").Value = -1351.75
$ws.Range("H104").Value = 1194
$ws.Range("J104").Value = 1950
$ws.Range("L104").Value = 5850
$ws.Range("N104").Value = -9344
$ws.Range("H113").Value = 4664.923
$ws.Range("I113").Value = 5736.852
$ws.Range("J113").Value = 2253.0833
$ws.Range("K113").Value = 5736.852
$ws.Range("L113").Value = 2253.0833
$ws.Range("M113").Value = -2482.852
$ws.Range("N113").Value = -8761.0833
$ws.Range("H122").Value = 3343.4285
$ws.Range("I122").Value = 2849.75
$ws.Range("K122").Value = 8549.25
$ws.Range("M122").Value = -6099.25
$ws.Range("H129").Value = 1646.3334
$ws.Range("I129").Value = 1242.1428
$ws.Range("K129").Value = 3726.4284
$ws.Range("M129").Value = 1273.5716
$ws.Range("H132").Value = 2170.7727
$ws.Range("I132").Value = 2170.7727
$ws.Range("K132").Value = 6512.3181
$ws.Range("M132").Value = -3982.3181
$ws.Range("H135").Value = 2504
$ws.Range("I135").Value = 2504
$ws.Range("K135").Value = 22536
$ws.Range("M135").Value = -20001
$ws.Range("H137").Value = 3032.95
$ws.Range("I137").Value = 2578.6365
$ws.Range("J137").Value = 3588.2222
$ws.Range("K137").Value = 7735.9095
$ws.Range("L137").Value = 10764.6666
$ws.Range("M137").Value = -5185.9095
$ws.Range("N137").Value = -15864.6666
$ws.Range("H138").Value = 1882.1111
$ws.Range("J138").Value = 2460.1587
$ws.Range("L138").Value = 7380.4761
$ws.Range("N138").Value = -17660.4761

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2829.9167
$ws.Range("I2").Value = 2551.2222
$ws.Range("K2").Value = 2551.2222
$ws.Range("M2").Value = -2438.2222
$ws.Range("H45").Value = 2163.75
$ws.Range("I45").Value = 1385
$ws.Range("K45").Value = 1385
$ws.Range("M45").Value = -1008
$ws.Range("H63").Value = 4072.8096
$ws.Range("I63").Value = 2328.5
$ws.Range("K63").Value = 2328.5
$ws.Range("M63").Value = -1642.5
$ws.Range("H66").Value = 4072.8096
$ws.Range("I66").Value = 2328.5
$ws.Range("K66").Value = 11642.5
$ws.Range("M66").Value = -8210.5
$ws.Range("H74").Value = 2887.2163
$ws.Range("I74").Value = 2416.5312
$ws.Range("K74").Value = 2416.5312
$ws.Range("M74").Value = -1542.5312
$ws.Range("H77").Value = 2887.2163
$ws.Range("I77").Value = 2416.5312
$ws.Range("K77").Value = 12082.656
$ws.Range("M77").Value = -7714.655999999999
$ws.Range("H110").Value = 3347
$ws.Range("I110").Value = 3194.5
$ws.Range("J110").Value = 3499.5
$ws.Range("K110").Value = 3194.5
$ws.Range("L110").Value = 3499.5
$ws.Range("M110").Value = -1149.5
$ws.Range("N110").Value = -7589.5
$ws.Range("H116").Value = 2829.9167
$ws.Range("I116").Value = 2551.2222
$ws.Range("K116").Value = 2551.2222
$ws.Range("M116").Value = -257.2222000000002
$ws.Range("H122").Value = 4018.3333
$ws.Range("I122").Value = 3531.5
$ws.Range("J122").Value = 5965.6665
$ws.Range("K122").Value = 10594.5
$ws.Range("L122").Value = 17896.9995
$ws.Range("M122").Value = -8144.5
$ws.Range("N122").Value = -22796.9995
$ws.Range("H132").Value = 420809.84
$ws.Range("I132").Value = 629064.25
$ws.Range("K132").Value = 1887192.75
$ws.Range("M132").Value = -1884662.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2829.9167
$ws.Range("I3").Value = 2551.2222
$ws.Range("K3").Value = 2551.2222
$ws.Range("M3").Value = -2437.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1779
$ws.Range("I10").Value = 1667.3334
$ws.Range("K10").Value = 1667.3334
$ws.Range("M10").Value = -1528.3334
$ws.Range("H31").Value = 4519.3335
$ws.Range("I31").Value = 2286.5625
$ws.Range("J31").Value = 6824.129
$ws.Range("K31").Value = 2286.5625
$ws.Range("L31").Value = 6824.129
$ws.Range("M31").Value = -1991.5625
$ws.Range("N31").Value = -7414.129
$ws.Range("H34").Value = 4519.3335
$ws.Range("I34").Value = 2286.5625
$ws.Range("J34").Value = 6824.129
$ws.Range("K34").Value = 2286.5625
$ws.Range("L34").Value = 6824.129
$ws.Range("M34").Value = -2084.5625
$ws.Range("N34").Value = -7228.129
$ws.Range("H99").Value = 3167.375
$ws.Range("I99").Value = 3048.4285
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 3048.4285
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -1550.4285
$ws.Range("N99").Value = -6996
$ws.Range("H126").Value = 3167.375
$ws.Range("I126").Value = 3048.4285
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 9145.2855
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -6675.2855
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 7757.091
$ws.Range("J137").Value = 9134.223
$ws.Range("L137").Value = 27402.669
$ws.Range("N137").Value = -37602.669

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1372.75
$ws.Range("I122").Value = 1496.6666
$ws.Range("K122").Value = 4489.9998
$ws.Range("M122").Value = -2039.9998
$ws.Range("H126").Value = 2777.875
$ws.Range("I126").Value = 2423.75
$ws.Range("K126").Value = 7271.25
$ws.Range("M126").Value = -4801.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3725.8572
$ws.Range("I7").Value = 3680.1667
$ws.Range("K7").Value = 3680.1667
$ws.Range("M7").Value = -3568.1667
$ws.Range("H16").Value = 286.875
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1340
$ws.Range("H126").Value = 3725.8572
$ws.Range("I126").Value = 3680.1667
$ws.Range("K126").Value = 11040.5001
$ws.Range("M126").Value = -8570.500100000001
$ws.Range("H132").Value = 719500
$ws.Range("I132").Value = 1004500.1
$ws.Range("K132").Value = 3013500.3
$ws.Range("M132").Value = -3010970.3
$ws.Range("H136").Value = 2521.0527
$ws.Range("I136").Value = 2457.0715
$ws.Range("K136").Value = 7371.2145
$ws.Range("M136").Value = -4821.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5505.7144
$ws.Range("I62").Value = 4280
$ws.Range("J62").Value = 6425
$ws.Range("K62").Value = 4280
$ws.Range("L62").Value = 6425
$ws.Range("M62").Value = -3656
$ws.Range("N62").Value = -7673
$ws.Range("H65").Value = 5505.7144
$ws.Range("I65").Value = 4280
$ws.Range("J65").Value = 6425
$ws.Range("K65").Value = 21400
$ws.Range("L65").Value = 32125
$ws.Range("M65").Value = -18280
$ws.Range("N65").Value = -38365
$ws.Range("H96").Value = 13963.833
$ws.Range("I96").Value = 7874.25
$ws.Range("K96").Value = 7874.25
$ws.Range("M96").Value = -6501.25
$ws.Range("H132").Value = 32756.06
$ws.Range("I132").Value = 43686.457
$ws.Range("J132").Value = 3608.3333
$ws.Range("K132").Value = 131059.371
$ws.Range("L132").Value = 10824.9999
$ws.Range("M132").Value = -128529.371
$ws.Range("N132").Value = -15884.9999
$ws.Range("H136").Value = 2077
$ws.Range("I136").Value = 1554.9166
$ws.Range("J136").Value = 3330
$ws.Range("K136").Value = 4664.7498
$ws.Range("L136").Value = 9990
$ws.Range("M136").Value = -2114.7498
$ws.Range("N136").Value = -15090
